$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "Supplied" header in column N, row 1, matching the style
# used by the rest of the header row (s="1").
$ws.Range("N1").Value = "Supplied"

# Reuse the exact same cell format as the other header cells (e.g. M1)
# by copying formats only, so no new style entries are created.
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Move the active selection to N2, as recorded in the saved workbook.
$ws.Range("N2").Select()
